$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.390.85"
$ws.Range("E2").Value = "  -4.38%  "
$ws.Range("D3").Value = "1.569.20"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "290.90"
$ws.Range("E6").Value = "  -2.73%  "
$ws.Range("D7").Value = "0.3683"
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("D8").Value = "49.55"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").Value = "0.3398"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").Value = "1.174"
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").Value = "0.07601"
$ws.Range("E11").Value = "  -5.86%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "21.20"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").Value = "6.049"
$ws.Range("D15").Value = "6.903"
$ws.Range("E15").Value = "  -5.36%  "
$ws.Range("E16").Value = "  -5.13%  "
$ws.Range("D17").Value = "1.573.84"
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("D18").Value = "89.11"
$ws.Range("E18").Value = "  -7.90%  "
$ws.Range("D19").Value = "0.06779"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "6.236"
$ws.Range("E21").Value = "  -7.10%  "
$ws.Range("D22").Value = "0.5359"
$ws.Range("E22").Value = "  -6.27%  "
$ws.Range("D23").Value = "16.51"
$ws.Range("E23").Value = "  -4.79%  "
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").Value = "22.410.32"
$ws.Range("E25").Value = "  -4.34%  "
$ws.Range("D26").Value = "2.385"
$ws.Range("E26").Value = "  -3.75%  "
$ws.Range("D27").Value = "2.984"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").Value = "19.88"
$ws.Range("E28").Value = "  -4.66%  "
$ws.Range("D29").Value = "145.61"
$ws.Range("E29").Value = "  -4.48%  "
$ws.Range("D30").Value = "4.968"
$ws.Range("E30").Value = "  -4.71%  "
$ws.Range("D31").Value = "125.54"
$ws.Range("E31").Value = "  -5.27%  "
$ws.Range("D32").Value = "1.757.10"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").Value = "1.046"
$ws.Range("E33").Value = "  +6.55%  "
$ws.Range("D34").Value = "6.256"
$ws.Range("E34").Value = "  -8.69%  "
$ws.Range("D35").Value = "1.995"
$ws.Range("E35").Value = "  -6.30%  "
$ws.Range("D36").Value = "10.31"
$ws.Range("E36").Value = "  -10.19%  "
$ws.Range("D37").Value = "0.08459"
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("D38").Value = "0.02546"
$ws.Range("E38").Value = "  -5.76%  "
$ws.Range("D39").Value = "0.2328"
$ws.Range("D40").Value = "0.06558"
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("D41").Value = "5.533"
$ws.Range("E41").Value = "  -6.31%  "
$ws.Range("E42").Value = "  -7.86%  "
$ws.Range("D43").Value = "1.247"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "0.6369"
$ws.Range("E44").Value = "  -7.02%  "
$ws.Range("D45").Value = "14.34"
$ws.Range("E45").Value = "  -7.42%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "0.5992"
$ws.Range("E47").Value = "  -5.18%  "
$ws.Range("D48").Value = "3.774"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("D49").Value = "2.136"
$ws.Range("E49").Value = "  -5.02%  "
$ws.Range("D50").Value = "1.260"
$ws.Range("E50").Value = "  +7.22%  "
$ws.Range("D51").Value = "123.57"
$ws.Range("E51").Value = "  -2.61%  "
